$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (old rows 3-6 shift down to 4-7)
$ws.Rows.Item(3).Insert()

# The new row 3's K3:L3 should carry the same wrap-text/shaded style as K2:L2
$ws.Range("K2:L2").Copy()
$ws.Range("K3:L3").PasteSpecial(-4122)

# --- New row 3: "Multi test" entry ---
$ws.Range("A3").Value = 2010
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 30
$ws.Range("K3").Value = "Multi test"
$ws.Range("L3").Value = "Test for multi"
$ws.Range("M3").Value = "exampleMulti"

# --- Row 4 (was row 3): Media column now references the example asset ---
$ws.Range("M4").Value = "exampleImg"

# --- Row 5 (was row 4): Media column now references the example asset, plus new Q5 ---
$ws.Range("M5").Value = "exampleAud"
$ws.Range("Q5").Value = "exampleImg"

# --- Row 6 (was row 5): Media column now references the example asset ---
$ws.Range("M6").Value = "exampleVid"

# --- Row 7 (was row 6): Media gains a value, Q7 media caption changes ---
$ws.Range("M7").Value = "demo2"
$ws.Range("Q7").Value = "exampleimg"

# --- New row 8 ---
$ws.Range("A8").Value = 2010
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 21
$ws.Range("K8").Value = "new image set"
$ws.Range("M8").Value = "demo"

# --- New row 9 ---
$ws.Range("A9").Value = 2010
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 22
$ws.Range("K9").Value = "new image set"
$ws.Range("L9").Value = "Here will be 5 images"
$ws.Range("M9").Value = "demo2"

# --- New row 10 ---
$ws.Range("A10").Value = 2005
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 19
$ws.Range("K10").Value = "out of order"

# Column width tweaks: P:Q become one uniform (wider) width, losing their bestFit auto-sizing
$ws.Range("P1:Q1").EntireColumn.ColumnWidth = 16.6

# View: zoom to 85%, and select K10 (last edited cell) as the active cell
$win = $excel.ActiveWindow
$win.Zoom = 85
$ws.Range("K10").Select()
